# Adding NUnit Test Case for Excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: change the Test Action from "Verify WebElement Availability" to
# "ClickElementByXPath" and clear the now-unused D6 value ("write").
$ws.Range("B6").Value = "ClickElementByXPath"
$ws.Range("D6").ClearContents()

# Row 7: the whole test step row is removed, leaving only the styled
# (italic) empty cell at B7.
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()

# Update the active selection to D5:D6 with D5 as the active cell.
$ws.Range("D5:D6").Select()
